$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 4098.6
$ws.Range("I13").Value = 995
$ws.Range("J13").Value = 4874.5
$ws.Range("K13").Value = 995
$ws.Range("L13").Value = 4874.5
$ws.Range("M13").Value = -826
$ws.Range("N13").Value = -5212.5

$ws.Range("H40").Value = 1131.4073
$ws.Range("I40").Value = 1041.1765
$ws.Range("J40").Value = 1284.8
$ws.Range("K40").Value = 1041.1765
$ws.Range("L40").Value = 1284.8
$ws.Range("M40").Value = -866.1765
$ws.Range("N40").Value = -1634.8

$ws.Range("H53").Value = 272.23077
$ws.Range("I53").Value = 181.33333
$ws.Range("J53").Value = 350.14285
$ws.Range("K53").Value = 181.33333
$ws.Range("L53").Value = 350.14285
$ws.Range("M53").Value = 455.66667
$ws.Range("N53").Value = -1624.14285

$ws.Range("H70").Value = 3517.5293
$ws.Range("I70").Value = 2250
$ws.Range("J70").Value = 3686.5334
$ws.Range("K70").Value = 6750
$ws.Range("L70").Value = 11059.6002
$ws.Range("M70").Value = -6480
$ws.Range("N70").Value = -11599.6002

$ws.Range("H73").Value = 3517.5293
$ws.Range("I73").Value = 2250
$ws.Range("J73").Value = 3686.5334
$ws.Range("K73").Value = 6750
$ws.Range("L73").Value = 11059.6002
$ws.Range("M73").Value = -5814
$ws.Range("N73").Value = -12931.6002

$ws.Range("H80").Value = 977.29266
$ws.Range("I80").Value = 1013.5294
$ws.Range("J80").Value = 951.625
$ws.Range("K80").Value = 3040.5882
$ws.Range("L80").Value = 2854.875
$ws.Range("M80").Value = -2042.5882
$ws.Range("N80").Value = -4850.875

$ws.Range("H83").Value = 977.29266
$ws.Range("I83").Value = 1013.5294
$ws.Range("J83").Value = 951.625
$ws.Range("K83").Value = 9121.7646
$ws.Range("L83").Value = 8564.625
$ws.Range("M83").Value = -4129.7646
$ws.Range("N83").Value = -18548.625

$ws.Range("H92").Value = 5273.905
$ws.Range("I92").Value = 5512.6
$ws.Range("K92").Value = 5512.6
$ws.Range("M92").Value = -4264.6

$ws.Range("H132").Value = 75428.37
$ws.Range("I132").Value = 41434.64
$ws.Range("K132").Value = 124303.92
$ws.Range("M132").Value = -121773.92

$ws.Range("H136").Value = 85579.5
$ws.Range("J136").Value = 85579.5
$ws.Range("L136").Value = 85579.5
$ws.Range("N136").Value = -95779.5

$ws.Range("H139").Value = 152484.75
$ws.Range("J139").Value = 184969.5
$ws.Range("L139").Value = 184969.5
$ws.Range("N139").Value = -195249.5

$ws.Range("H141").Value = 2201.7
$ws.Range("I141").Value = 1002.4286
$ws.Range("K141").Value = 3007.2858
$ws.Range("M141").Value = 2172.7142

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H12").Value = 23600.5
$ws.Range("I12").Value = 701.5
$ws.Range("K12").Value = 701.5
$ws.Range("M12").Value = -528.5

$ws.Range("H14").Value = 4278.3335
$ws.Range("I14").Value = 3735.3333
$ws.Range("J14").Value = 4821.3335
$ws.Range("K14").Value = 3735.3333
$ws.Range("L14").Value = 4821.3335
$ws.Range("M14").Value = -3560.3333
$ws.Range("N14").Value = -5171.3335

$ws.Range("H32").Value = 3105.8281
$ws.Range("I32").Value = 2490.625
$ws.Range("K32").Value = 2490.625
$ws.Range("M32").Value = -2203.625

$ws.Range("H97").Value = 985.4231
$ws.Range("I97").Value = 964.84
$ws.Range("K97").Value = 964.84
$ws.Range("M97").Value = -468.84

$ws.Range("H122").Value = 1470.3529
$ws.Range("I122").Value = 880.8570999999999
$ws.Range("K122").Value = 2642.5713
$ws.Range("M122").Value = -192.5712999999996

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 1332.2858
$ws.Range("J80").Value = 1426.091
$ws.Range("L80").Value = 1426.091
$ws.Range("N80").Value = -3422.091

$ws.Range("H83").Value = 1332.2858
$ws.Range("J83").Value = 1426.091
$ws.Range("L83").Value = 7130.455
$ws.Range("N83").Value = -17114.455

$ws.Range("H88").Value = 12541438
$ws.Range("J88").Value = 14330214
$ws.Range("L88").Value = 14330214
$ws.Range("N88").Value = -14331026

$ws.Range("H91").Value = 12541438
$ws.Range("J91").Value = 14330214
$ws.Range("L91").Value = 14330214
$ws.Range("N91").Value = -14333022

$ws.Range("H94").Value = 1727.7778
$ws.Range("I94").Value = 1312.9032
$ws.Range("K94").Value = 1312.9032
$ws.Range("M94").Value = -861.9032

$ws.Range("H107").Value = 1160.4193
$ws.Range("I107").Value = 1268.05
$ws.Range("J107").Value = 964.7273
$ws.Range("K107").Value = 1268.05
$ws.Range("L107").Value = 964.7273
$ws.Range("M107").Value = 651.95
$ws.Range("N107").Value = -4804.7273

$ws.Range("H132").Value = 85166.664
$ws.Range("J132").Value = 85166.664
$ws.Range("L132").Value = 85166.664
$ws.Range("N132").Value = -95286.664

$ws.Range("H133").Value = 99968
$ws.Range("J133").Value = 99968
$ws.Range("L133").Value = 99968
$ws.Range("N133").Value = -110088

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 415.9091
$ws.Range("I12").Value = 655.1667
$ws.Range("J12").Value = 128.8
$ws.Range("K12").Value = 655.1667
$ws.Range("L12").Value = 128.8
$ws.Range("M12").Value = -485.1667
$ws.Range("N12").Value = -468.8

$ws.Range("H58").Value = 1414.9231
$ws.Range("I58").Value = 1414.9231
$ws.Range("K58").Value = 1414.9231
$ws.Range("M58").Value = -1211.9231

$ws.Range("H136").Value = 1414.9231
$ws.Range("I136").Value = 1414.9231
$ws.Range("K136").Value = 4244.7693
$ws.Range("M136").Value = -1694.7693

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H19").Value = 300
$ws.Range("J19").Value = 500
$ws.Range("L19").Value = 1500
$ws.Range("N19").Value = -1848

$ws.Range("H23").Value = 213.42857
$ws.Range("I23").Value = 146.38461
$ws.Range("K23").Value = 439.15383
$ws.Range("M23").Value = -204.15383

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 1099.5
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()

$ws.Range("H38").Value = 41012
$ws.Range("I38").Value = 60000
$ws.Range("J38").Value = 22024
$ws.Range("K38").Value = 60000
$ws.Range("L38").Value = 22024
$ws.Range("M38").Value = -59537
$ws.Range("N38").Value = -22950

$ws.Range("H97").Value = 1117.9445
$ws.Range("I97").Value = 867.3570999999999
$ws.Range("K97").Value = 867.3570999999999
$ws.Range("M97").Value = -371.3570999999999

$ws.Range("H122").Value = 2633.762
$ws.Range("I122").Value = 2522.5
$ws.Range("J122").Value = 3301.3333
$ws.Range("K122").Value = 7567.5
$ws.Range("L122").Value = 9903.999899999999
$ws.Range("M122").Value = -5117.5
$ws.Range("N122").Value = -14803.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2071
$ws.Range("I82").Value = 1805.909
$ws.Range("J82").Value = 2800
$ws.Range("K82").Value = 1805.909
$ws.Range("L82").Value = 2800
$ws.Range("M82").Value = -1444.909
$ws.Range("N82").Value = -3522

$ws.Range("H85").Value = 2071
$ws.Range("I85").Value = 1805.909
$ws.Range("J85").Value = 2800
$ws.Range("K85").Value = 1805.909
$ws.Range("L85").Value = 2800
$ws.Range("M85").Value = -557.9090000000001
$ws.Range("N85").Value = -5296

$ws.Range("H132").Value = 5418.727
$ws.Range("I132").Value = 4075.75
$ws.Range("K132").Value = 12227.25
$ws.Range("M132").Value = -9697.25

$ws.Range("H136").Value = 3969.1428
$ws.Range("I136").Value = 3745.9565
$ws.Range("K136").Value = 11237.8695
$ws.Range("M136").Value = -8687.869499999999

$ws.Range("H141").Value = 89250
$ws.Range("J141").Value = 89250
$ws.Range("L141").Value = 89250
$ws.Range("N141").Value = -99610

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H34").Value = 22200
$ws.Range("J34").Value = 22200
$ws.Range("L34").Value = 22200
$ws.Range("N34").Value = -22606

$ws.Range("H42").Value = 23200
$ws.Range("J42").Value = 23200
$ws.Range("L42").Value = 23200
$ws.Range("N42").Value = -23956

$ws.Range("H81").Value = 9858
$ws.Range("I81").Value = 27111
$ws.Range("K81").Value = 54222
$ws.Range("M81").Value = -53161

$ws.Range("H84").Value = 9858
$ws.Range("I84").Value = 27111
$ws.Range("K84").Value = 271110
$ws.Range("M84").Value = -265806

$ws.Range("H107").Value = 1280.2778
$ws.Range("I107").Value = 722.9
$ws.Range("K107").Value = 2168.7
$ws.Range("M107").Value = -248.6999999999998

$ws.Range("H126").Value = 4901.3
$ws.Range("I126").Value = 4668.1113
$ws.Range("K126").Value = 14004.3339
$ws.Range("M126").Value = -11534.3339

$ws.Range("H137").Value = 87659
$ws.Range("J137").Value = 83244
$ws.Range("L137").Value = 83244
$ws.Range("N137").Value = -93444

$ws.Range("H139").Value = 105575
$ws.Range("I139").Value = 120650
$ws.Range("J139").Value = 90500
$ws.Range("K139").Value = 120650
$ws.Range("L139").Value = 90500
$ws.Range("M139").Value = -115510
$ws.Range("N139").Value = -100780

$ws.Range("H141").Value = 110959.8
$ws.Range("J141").Value = 110959.8
$ws.Range("L141").Value = 110959.8
$ws.Range("N141").Value = -121319.8
